# Auto-generated edit script
# Applies scheduled market-data value updates to multiple sheets/cells
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 7876.9546
$ws.Range("I137").Value = 9664.294
$ws.Range("J137").Value = 1800
$ws.Range("K137").Value = 28992.882
$ws.Range("L137").Value = 5400
$ws.Range("M137").Value = -26442.882
$ws.Range("N137").Value = -10500

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2834.2693
$ws.Range("I110").Value = 2903.7917
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 2903.7917
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = -858.7917000000002
$ws.Range("N110").Value = -6090
$ws.Range("H132").Value = 23412.389
$ws.Range("I132").Value = 38761.855
$ws.Range("J132").Value = 2946.4285
$ws.Range("K132").Value = 116285.565
$ws.Range("L132").Value = 8839.2855
$ws.Range("M132").Value = -113755.565
$ws.Range("N132").Value = -13899.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1638.15
$ws.Range("I20").Value = 1260.9
$ws.Range("J20").Value = 2015.4
$ws.Range("K20").Value = 1260.9
$ws.Range("L20").Value = 2015.4
$ws.Range("M20").Value = -1013.9
$ws.Range("N20").Value = -2509.4
$ws.Range("H99").Value = 7022.222
$ws.Range("I99").Value = 8557.143
$ws.Range("J99").Value = 1650
$ws.Range("K99").Value = 8557.143
$ws.Range("L99").Value = 1650
$ws.Range("M99").Value = -7059.143
$ws.Range("N99").Value = -4646
$ws.Range("H134").Value = 4094.4048
$ws.Range("I134").Value = 4153.9697
$ws.Range("J134").Value = 3876
$ws.Range("K134").Value = 12461.9091
$ws.Range("L134").Value = 11628
$ws.Range("M134").Value = -9926.909099999999
$ws.Range("N134").Value = -16698

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 795.86664
$ws.Range("I16").Value = 853.3333
$ws.Range("J16").Value = 757.55554
$ws.Range("K16").Value = 853.3333
$ws.Range("L16").Value = 757.55554
$ws.Range("M16").Value = -566.3333
$ws.Range("N16").Value = -1331.55554
$ws.Range("H31").Value = 2977.1428
$ws.Range("I31").Value = 1204.5454
$ws.Range("J31").Value = 5976.923
$ws.Range("K31").Value = 1204.5454
$ws.Range("L31").Value = 5976.923
$ws.Range("M31").Value = -909.5454
$ws.Range("N31").Value = -6566.923
$ws.Range("H34").Value = 2977.1428
$ws.Range("I34").Value = 1204.5454
$ws.Range("J34").Value = 5976.923
$ws.Range("K34").Value = 1204.5454
$ws.Range("L34").Value = 5976.923
$ws.Range("M34").Value = -1002.5454
$ws.Range("N34").Value = -6380.923
$ws.Range("H58").Value = 2744.077
$ws.Range("I58").Value = 3106.077
$ws.Range("J58").Value = 2201.077
$ws.Range("K58").Value = 3106.077
$ws.Range("L58").Value = 2201.077
$ws.Range("M58").Value = -2903.077
$ws.Range("N58").Value = -2607.077
$ws.Range("H94").Value = 4647.7827
$ws.Range("I94").Value = 995.9231
$ws.Range("J94").Value = 9395.200000000001
$ws.Range("K94").Value = 995.9231
$ws.Range("L94").Value = 9395.200000000001
$ws.Range("M94").Value = -544.9231
$ws.Range("N94").Value = -10297.2
$ws.Range("H99").Value = 169352
$ws.Range("I99").Value = 335037.34
$ws.Range("J99").Value = 3666.6667
$ws.Range("K99").Value = 335037.34
$ws.Range("L99").Value = 3666.6667
$ws.Range("M99").Value = -333539.34
$ws.Range("N99").Value = -6662.6667
$ws.Range("H113").Value = 795.86664
$ws.Range("I113").Value = 853.3333
$ws.Range("J113").Value = 757.55554
$ws.Range("K113").Value = 853.3333
$ws.Range("L113").Value = 757.55554
$ws.Range("M113").Value = 1316.6667
$ws.Range("N113").Value = -5097.55554
$ws.Range("H126").Value = 169352
$ws.Range("I126").Value = 335037.34
$ws.Range("J126").Value = 3666.6667
$ws.Range("K126").Value = 1005112.02
$ws.Range("L126").Value = 11000.0001
$ws.Range("M126").Value = -1002642.02
$ws.Range("N126").Value = -15940.0001
$ws.Range("H132").Value = 2604.5
$ws.Range("I132").Value = 1061.3334
$ws.Range("J132").Value = 4147.6665
$ws.Range("K132").Value = 3184.0002
$ws.Range("L132").Value = 12442.9995
$ws.Range("M132").Value = -654.0001999999999
$ws.Range("N132").Value = -17502.9995
$ws.Range("H136").Value = 2744.077
$ws.Range("I136").Value = 3106.077
$ws.Range("J136").Value = 2201.077
$ws.Range("K136").Value = 9318.231
$ws.Range("L136").Value = 6603.231000000001
$ws.Range("M136").Value = -6768.231
$ws.Range("N136").Value = -11703.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 708.3333
$ws.Range("I80").Value = 483.33334
$ws.Range("K80").Value = 1450.00002
$ws.Range("M80").Value = -514.0000199999999
$ws.Range("H83").Value = 708.3333
$ws.Range("I83").Value = 483.33334
$ws.Range("K83").Value = 4350.00006
$ws.Range("M83").Value = 329.9999399999997
$ws.Range("H86").Value = 599.125
$ws.Range("H89").Value = 599.125
$ws.Range("H130").Value = 127601.625
$ws.Range("I130").Value = 780
$ws.Range("J130").Value = 145719
$ws.Range("K130").Value = 2340
$ws.Range("L130").Value = 437157
$ws.Range("M130").Value = 2680
$ws.Range("N130").Value = -447197
$ws.Range("H131").Value = 1242.1786
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 1299.2693
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 3897.8079
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -13977.8079
$ws.Range("H132").Value = 3950.3667
$ws.Range("J132").Value = 5214.0586
$ws.Range("L132").Value = 46926.52740000001
$ws.Range("N132").Value = -51986.52740000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4500
$ws.Range("I102").Value = 1333.3334
$ws.Range("J102").Value = 14000
$ws.Range("K102").Value = 1333.3334
$ws.Range("L102").Value = 14000
$ws.Range("M102").Value = 288.6666
$ws.Range("N102").Value = -17244
$ws.Range("H126").Value = 2305.25
$ws.Range("I126").Value = 2365.7144
$ws.Range("J126").Value = 2258.2222
$ws.Range("K126").Value = 7097.1432
$ws.Range("L126").Value = 6774.6666
$ws.Range("M126").Value = -4627.1432
$ws.Range("N126").Value = -11714.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 22729028
$ws.Range("I132").Value = 30304054
$ws.Range("J132").Value = 3954.2727
$ws.Range("K132").Value = 90912162
$ws.Range("L132").Value = 11862.8181
$ws.Range("M132").Value = -90909632
$ws.Range("N132").Value = -16922.8181
$ws.Range("H136").Value = 10644912
$ws.Range("I136").Value = 17878096
$ws.Range("J136").Value = 258801.88
$ws.Range("K136").Value = 53634288
$ws.Range("L136").Value = 776405.64
$ws.Range("M136").Value = -53631738
$ws.Range("N136").Value = -781505.64
